$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.367.78"

$ws.Range("D3").Value = "1.845.16"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'238.88"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").Value = "'0.6291"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.07553"
$ws.Range("E8").Value = "  -1.10%  "

$ws.Range("E9").Value = "  -1.03%  "

$ws.Range("D10").Value = "'24.59"
$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("D11").Value = "'0.07697"
$ws.Range("E11").Value = "  -0.32%  "

$ws.Range("D12").Value = "1.842.83"
$ws.Range("E12").Value = "  -5.27%  "

$ws.Range("D13").Value = "'4.979"

$ws.Range("E14").Value = "  -1.59%  "

$ws.Range("D15").Value = "'0.00001021"
$ws.Range("E15").Value = "  +2.39%  "

$ws.Range("D16").Value = "'83.05"
$ws.Range("E16").Value = "  -0.26%  "

$ws.Range("D17").Value = "2.091.06"
$ws.Range("E17").Value = "  -4.78%  "

$ws.Range("D18").Value = "'6.127"
$ws.Range("E18").Value = "  -1.05%  "

$ws.Range("D19").Value = "29.407.87"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'228.17"
$ws.Range("E20").Value = "  -2.10%  "

$ws.Range("E21").Value = "  -1.18%  "

$ws.Range("D22").Value = "'0.9998"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'7.451"
$ws.Range("E23").Value = "  -3.21%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "'156.71"
$ws.Range("E25").Value = "  +1.11%  "

$ws.Range("E26").Value = "  -0.21%  "

$ws.Range("D27").Value = "'8.350"
$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("D28").Value = "'17.62"
$ws.Range("E28").Value = "  -0.34%  "

$ws.Range("E29").Value = "  -1.39%  "

$ws.Range("E30").Value = "  +0.96%  "

$ws.Range("D31").Value = "'0.05628"
$ws.Range("E31").Value = "  -2.75%  "

$ws.Range("E32").Value = "  -0.29%  "

$ws.Range("D33").Value = "'4.038"
$ws.Range("E33").Value = "  +0.34%  "

$ws.Range("E34").Value = "  -2.52%  "

$ws.Range("D35").Value = "'1.155"
$ws.Range("E35").Value = "  -0.71%  "

$ws.Range("D36").Value = "'0.7157"
$ws.Range("E36").Value = "  -0.60%  "

$ws.Range("D37").Value = "'2.595"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").Value = "1.239.43"
$ws.Range("E38").Value = "  -0.66%  "

$ws.Range("D39").Value = "'0.01809"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("D40").Value = "'2.773"
$ws.Range("E40").Value = "  -0.78%  "

$ws.Range("D41").Value = "'6.217"
$ws.Range("E41").Value = "  +2.04%  "

$ws.Range("D42").Value = "'0.9014"
$ws.Range("E42").Value = "  -0.82%  "

$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").Value = "'101.80"
$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("D45").Value = "'65.76"
$ws.Range("E45").Value = "  -3.17%  "

$ws.Range("D46").Value = "'7.105"
$ws.Range("E46").Value = "  -2.58%  "

$ws.Range("E47").Value = "  -2.08%  "

$ws.Range("D48").Value = "'0.3996"
$ws.Range("E48").Value = "  -0.88%  "

$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'1.681"
$ws.Range("E49").Value = "  -1.19%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.952"
$ws.Range("E50").Value = "  -2.69%  "

$ws.Range("D51").Value = "'0.1118"
$ws.Range("E51").Value = "  -0.33%  "

